$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting existing data down
$ws.Rows.Item(1).Insert()

# Set the new header cell
$ws.Range("A1").Value = "ENSEMBL_ID"
